$wb = $excel.ActiveWorkbook

$wsNoise = $wb.Worksheets.Item("withNoise")
$wsNoise.Range("B3").Value = 346
$wsNoise.Range("C3").Value = 345
$wsNoise.Range("F3").Value = 350
$wsNoise.Range("G3").Value = 350
$wsNoise.Range("J3").Value = 349
$wsNoise.Range("K3").Value = 348
$wsNoise.Range("N3").Value = 344
$wsNoise.Range("O3").Value = 344
$wsNoise.Range("B4").Value = 0.9957368007414592
$wsNoise.Range("C4").Value = 0.9957917657510439
$wsNoise.Range("F4").Value = 0.9988934151795835
$wsNoise.Range("G4").Value = 0.9988934151795835
$wsNoise.Range("J4").Value = 0.07456617805400845
$wsNoise.Range("K4").Value = 0.07455375731348657
$wsNoise.Range("N4").Value = 0.09473645879962589
$wsNoise.Range("O4").Value = 0.09473645879962589
$wsNoise.Range("B5").Value = 0.995778058834679
$wsNoise.Range("C5").Value = 0.995801635395351
$wsNoise.Range("F5").Value = 0.9988651753657154
$wsNoise.Range("G5").Value = 0.9988651753657154
$wsNoise.Range("J5").Value = 0.07402426525385286
$wsNoise.Range("K5").Value = 0.07389053541194934
$wsNoise.Range("N5").Value = 0.09218837922826839
$wsNoise.Range("O5").Value = 0.09218837922826839
$wsNoise.Range("B6").Value = 0.9955987448874793
$wsNoise.Range("C6").Value = 0.9956189769887136
$wsNoise.Range("F6").Value = 0.9988722416207938
$wsNoise.Range("G6").Value = 0.9988722416207938
$wsNoise.Range("J6").Value = 0.07466250097719018
$wsNoise.Range("K6").Value = 0.07460929568291363
$wsNoise.Range("N6").Value = 0.09063134454432165
$wsNoise.Range("O6").Value = 0.09063134454432165
$wsNoise.Range("B7").Value = 0.9949839145655118
$wsNoise.Range("C7").Value = 0.9951115357807993
$wsNoise.Range("F7").Value = 0.9988643044524392
$wsNoise.Range("G7").Value = 0.9988643044524392
$wsNoise.Range("J7").Value = 0.07305110945056126
$wsNoise.Range("K7").Value = 0.07342755578251181
$wsNoise.Range("N7").Value = 0.08863230649054776
$wsNoise.Range("O7").Value = 0.08863230649054776
$wsNoise.Range("B30").Value = 323
$wsNoise.Range("C30").Value = 322
$wsNoise.Range("J30").Value = 323
$wsNoise.Range("K30").Value = 323
$wsNoise.Range("N30").Value = 319
$wsNoise.Range("O30").Value = 319
$wsNoise.Range("B31").Value = 0.06866314750687252
$wsNoise.Range("C31").Value = 0.06826055285127608
$wsNoise.Range("J31").Value = 0.06536514820025308
$wsNoise.Range("K31").Value = 0.06536514820025308
$wsNoise.Range("N31").Value = 0.0739992422096668
$wsNoise.Range("O31").Value = 0.0739992422096668
$wsNoise.Range("B32").Value = 0.07408269929472644
$wsNoise.Range("C32").Value = 0.07355318604759598
$wsNoise.Range("J32").Value = 0.07042566275682982
$wsNoise.Range("K32").Value = 0.07042566275682982
$wsNoise.Range("N32").Value = 0.07901123518687858
$wsNoise.Range("O32").Value = 0.07901123518687858
$wsNoise.Range("B33").Value = 0.07699301677801429
$wsNoise.Range("C33").Value = 0.07674966343422662
$wsNoise.Range("J33").Value = 0.07412178191509555
$wsNoise.Range("K33").Value = 0.07412178191509555
$wsNoise.Range("N33").Value = 0.08354979139152663
$wsNoise.Range("O33").Value = 0.08354979139152663
$wsNoise.Range("B34").Value = 0.07643071499573506
$wsNoise.Range("C34").Value = 0.07666627330717168
$wsNoise.Range("J34").Value = 0.07477522465328801
$wsNoise.Range("K34").Value = 0.07477522465328801
$wsNoise.Range("N34").Value = 0.085980507366749
$wsNoise.Range("O34").Value = 0.085980507366749

$wsNoNoise = $wb.Worksheets.Item("withoutNoise")
$wsNoNoise.Range("B3").Value = 347
$wsNoNoise.Range("C3").Value = 347
$wsNoNoise.Range("F3").Value = 350
$wsNoNoise.Range("G3").Value = 350
$wsNoNoise.Range("K3").Value = 350
$wsNoNoise.Range("N3").Value = 346
$wsNoNoise.Range("O3").Value = 346
$wsNoNoise.Range("B4").Value = 1
$wsNoNoise.Range("C4").Value = 1
$wsNoNoise.Range("F4").Value = 1
$wsNoNoise.Range("G4").Value = 1
$wsNoNoise.Range("K4").Value = 0
$wsNoNoise.Range("B5").Value = 1
$wsNoNoise.Range("C5").Value = 1
$wsNoNoise.Range("F5").Value = 1
$wsNoNoise.Range("G5").Value = 1
$wsNoNoise.Range("K5").Value = 0
$wsNoNoise.Range("B6").Value = 1
$wsNoNoise.Range("C6").Value = 1
$wsNoNoise.Range("F6").Value = 1
$wsNoNoise.Range("G6").Value = 1
$wsNoNoise.Range("K6").Value = 0
$wsNoNoise.Range("B7").Value = 1
$wsNoNoise.Range("C7").Value = 1
$wsNoNoise.Range("F7").Value = 1
$wsNoNoise.Range("G7").Value = 1
$wsNoNoise.Range("K7").Value = 0
